$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A158").Value = "2023-12-10 13:01:41"
$ws.Range("B158").Value = 0.0004

$ws.Range("A159").Value = "2023-12-10 13:01:54"
$ws.Range("B159").Value = 0.0006000000000000001
